$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("samples_retained")

# Row 22 (oreau2) - fill in the remaining details
$ws.Range("B22").Value = "acted"
$ws.Range("C22").Value = 62
$ws.Range("D22").Value = 302
$ws.Range("E22").Value = 70
$ws.Range("H22").Value = 32
$ws.Range("I22").Value = "Joy, Sadness, Fear, Anger, Surprise, Disgust, Neutral"

# New dataset rows 23-28
$ws.Range("A23").Value = "ravdess"
$ws.Range("A24").Value = "savee"
$ws.Range("A25").Value = "ShEMO"
$ws.Range("A26").Value = "tess"
$ws.Range("A27").Value = "urdu"
$ws.Range("A28").Value = "vivae"

# Match the final selection state recorded in the workbook
[void]$ws.Range("B23").Select()

